$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. These D-column cells store prices/volumes
# as plain text (inline strings) rather than numbers, so we force the
# cell format to Text before assigning, ensuring values like "6.380" or
# "0.003350" keep their exact textual representation (no numeric coercion,
# no dropped trailing zeros).
$updates = @{
    "D2"  = "274.58"
    "D3"  = "23.04"
    "D4"  = "6.380"
    "D6"  = "3.659"
    "D7"  = "6.671"
    "D8"  = "1.358"
    "D11" = "0.1632"
    "D12" = "0.08323"
    "D13" = "0.03425"
    "D14" = "0.03081"
    "D15" = "0.09302"
    "D16" = "3.856"
    "D17" = "0.001641"
    "D18" = "0.04779"
    "D19" = "0.006409"
    "D21" = "0.001092"
    "D25" = "0.3345"
    "D40" = "0.04692"
    "D41" = "0.007051"
    "D42" = "0.1162"
    "D43" = "0.003350"
    "D44" = "0.01183"
    "D45" = "0.00006262"
    "D48" = "0.7962"
    "D49" = "0.04033"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
